# dsa dp and bits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 97: "Easay" -> "Easy" (typo fix on Difficulty column) ---
$ws.Range("B97").Value = "Easy"

$tbl = $ws.ListObjects.Item(1)

# --- Add row 98: 97. Interleaving String (Medium / Dynamic Programming) ---
$tbl.ListRows.Add() | Out-Null
$ws.Range("E98").Value = "https://leetcode.com/problems/interleaving-string/solutions/32078/dp-solution-in-java/ "
$ws.Hyperlinks.Add($ws.Range("E98"), "https://leetcode.com/problems/interleaving-string/solutions/32078/dp-solution-in-java/") | Out-Null
$ws.Range("A90:E90").Copy()
$ws.Range("A98:E98").PasteSpecial(-4122) | Out-Null
$ws.Range("A98").Value = "97. Interleaving String"
$ws.Range("B98").Value = "Medium"
$ws.Range("C98").Value = "Dynamic Programming"
$ws.Range("D98").Value = "The position in the target string s3 is given by the sum of the positions of s1 and s2. In the DP, out of bounds is the base case. We need 1 extra row and column, as it can be the case where we take the last character from 1 string, while the other has ended. Start from the bottom right, then work to the top left."
$ws.Range("E98").Value = "https://leetcode.com/problems/interleaving-string/solutions/32078/dp-solution-in-java/ "

# --- Add row 99: 329. Longest Increasing Path in a Matrix (Hard / Dynamic Programming) ---
$tbl.ListRows.Add() | Out-Null
$ws.Range("E99").Value = "https://leetcode.com/problems/longest-increasing-path-in-a-matrix/solutions/78308/15ms-concise-java-solution/ "
$ws.Hyperlinks.Add($ws.Range("E99"), "https://leetcode.com/problems/longest-increasing-path-in-a-matrix/solutions/78308/15ms-concise-java-solution/") | Out-Null
$ws.Range("A91:E91").Copy()
$ws.Range("A99:E99").PasteSpecial(-4122) | Out-Null
$ws.Range("A99").Value = "329. Longest Increasing Path in a Matrix"
$ws.Range("B99").Value = "Hard"
$ws.Range("C99").Value = "Dynamic Programming"
$ws.Range("D99").Value = "We use DFS on a grid. We cannot reuse positions as an implicit condition of the longest increasing path. We store the longest increasing path from each position in the grid, to cache repeated work. The crux of the solution is the state transition/subproblem calculation where we take 1 + the maximum of the DFS of all 4 directions, and place it in the position in the dp matrix. We encode into our boundary conditions check the condition that value is greater than its parent, to maintain an increasing path.  Key implementation details are: the initialization of the base state of -1, flagging the cell as unvisited. DFS on a 2D matrix 2nd pass to find the max value."
$ws.Range("E99").Value = "https://leetcode.com/problems/longest-increasing-path-in-a-matrix/solutions/78308/15ms-concise-java-solution/ "
$ws.Range("D99").WrapText = $true
$ws.Rows.Item(99).RowHeight = 15

# --- Add row 100: 268. Missing Number (Easy / Bit Manipulation) ---
$tbl.ListRows.Add() | Out-Null
$ws.Range("E100").Value = "https://leetcode.com/problems/missing-number/solutions/69791/4-line-simple-java-bit-manipulate-solution-with-explaination/ "
$ws.Hyperlinks.Add($ws.Range("E100"), "https://leetcode.com/problems/missing-number/solutions/69791/4-line-simple-java-bit-manipulate-solution-with-explaination/") | Out-Null
$ws.Range("A97:E97").Copy()
$ws.Range("A100:E100").PasteSpecial(-4122) | Out-Null
$ws.Range("A100").Value = "268. Missing Number"
$ws.Range("B100").Value = "Easy"
$ws.Range("C100").Value = "Bit Manipulation"
$ws.Range("D100").Value = "The non bit manipulation solution is to simply take the sum of input and target arrays, then take the difference as the missing number. We can use XOR for O(1) memory and O(N) time. E.g. a^b^b = a. If the numbers are the same, XOR will give a 0 in the output. The order of the numbers does not matter. We XOR the input and the target, and therefore only the missing number is remaining."
$ws.Range("E100").Value = "https://leetcode.com/problems/missing-number/solutions/69791/4-line-simple-java-bit-manipulate-solution-with-explaination/ "

# --- Update selection to match the final state ---
$ws.Range("D100").Select()

Write-Output "edit complete"
